$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.238.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.214.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.78%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.434"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.770.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.288.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.223.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.24%  "
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.803.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0708"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.255.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.811"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
